# Auto-generated script applying cryptos.xlsx price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "30.102.41"
Set-TextValue $ws.Range("E2") "  -1.56%  "
Set-TextValue $ws.Range("D3") "2.105.41"
Set-TextValue $ws.Range("E3") "  -0.26%  "
Set-TextValue $ws.Range("D4") "1.003"
Set-TextValue $ws.Range("E4") "  -0.93%  "
Set-TextValue $ws.Range("D5") "351.14"
Set-TextValue $ws.Range("E5") "  +4.44%  "
Set-TextValue $ws.Range("D6") "1.003"
Set-TextValue $ws.Range("E6") "  -0.76%  "
Set-TextValue $ws.Range("D7") "0.5170"
Set-TextValue $ws.Range("E7") "  -1.46%  "
Set-TextValue $ws.Range("D8") "0.4488"
Set-TextValue $ws.Range("E8") "  -1.38%  "
Set-TextValue $ws.Range("D9") "52.73"
Set-TextValue $ws.Range("E9") "  -4.14%  "
Set-TextValue $ws.Range("D10") "0.08963"
Set-TextValue $ws.Range("E10") "  -0.32%  "
Set-TextValue $ws.Range("D11") "1.174"
Set-TextValue $ws.Range("E11") "  +0.46%  "
Set-TextValue $ws.Range("D12") "25.65"
Set-TextValue $ws.Range("E12") "  +4.64%  "
Set-TextValue $ws.Range("D13") "2.103.19"
Set-TextValue $ws.Range("E13") "  -0.29%  "
Set-TextValue $ws.Range("D14") "6.759"
Set-TextValue $ws.Range("E14") "  -1.49%  "
Set-TextValue $ws.Range("D15") "8.144"
Set-TextValue $ws.Range("E15") "  +0.36%  "
Set-TextValue $ws.Range("D16") "99.56"
Set-TextValue $ws.Range("E16") "  +2.54%  "
Set-TextValue $ws.Range("D17") "0.00001150"
Set-TextValue $ws.Range("E17") "  -2.02%  "
Set-TextValue $ws.Range("D18") "1.004"
Set-TextValue $ws.Range("E18") "  -0.75%  "
Set-TextValue $ws.Range("D19") "20.56"
Set-TextValue $ws.Range("E19") "  +6.21%  "
Set-TextValue $ws.Range("E20") "  -0.04%  "
Set-TextValue $ws.Range("D21") "1.003"
Set-TextValue $ws.Range("E21") "  -0.72%  "
Set-TextValue $ws.Range("D22") "6.236"
Set-TextValue $ws.Range("E22") "  -0.33%  "
Set-TextValue $ws.Range("D23") "30.189.53"
Set-TextValue $ws.Range("E23") "  -1.50%  "
Set-TextValue $ws.Range("D24") "12.87"
Set-TextValue $ws.Range("E24") "  +0.38%  "
Set-TextValue $ws.Range("D25") "2.349"
Set-TextValue $ws.Range("E25") "  -0.51%  "
Set-TextValue $ws.Range("D26") "2.352.96"
Set-TextValue $ws.Range("E26") "  -0.25%  "
Set-TextValue $ws.Range("D27") "22.07"
Set-TextValue $ws.Range("E27") "  -1.15%  "
Set-TextValue $ws.Range("D28") "2.556"
Set-TextValue $ws.Range("E28") "  +1.09%  "
Set-TextValue $ws.Range("D29") "162.68"
Set-TextValue $ws.Range("E29") "  -0.48%  "
Set-TextValue $ws.Range("D30") "133.74"
Set-TextValue $ws.Range("E30") "  +0.04%  "
Set-TextValue $ws.Range("D31") "1.184"
Set-TextValue $ws.Range("E31") "  -3.26%  "
Set-TextValue $ws.Range("D32") "0.1067"
Set-TextValue $ws.Range("E32") "  -0.40%  "
Set-TextValue $ws.Range("D33") "1.653"
Set-TextValue $ws.Range("E33") "  +1.69%  "
Set-TextValue $ws.Range("D34") "6.265"
Set-TextValue $ws.Range("E34") "  -1.36%  "
Set-TextValue $ws.Range("D35") "3.962"
Set-TextValue $ws.Range("E35") "  -0.10%  "
Set-TextValue $ws.Range("D36") "5.937"
Set-TextValue $ws.Range("E36") "  +1.15%  "
Set-TextValue $ws.Range("D37") "10.19"
Set-TextValue $ws.Range("E37") "  -2.49%  "
Set-TextValue $ws.Range("D38") "0.02588"
Set-TextValue $ws.Range("E38") "  -0.25%  "
Set-TextValue $ws.Range("D39") "0.06841"
Set-TextValue $ws.Range("E39") "  +0.23%  "
Set-TextValue $ws.Range("D40") "0.2310"
Set-TextValue $ws.Range("E40") "  -0.54%  "
Set-TextValue $ws.Range("D41") "12.55"
Set-TextValue $ws.Range("E41") "  -0.54%  "
Set-TextValue $ws.Range("D42") "0.6834"
Set-TextValue $ws.Range("E42") "  -0.35%  "
Set-TextValue $ws.Range("D43") "1.256"
Set-TextValue $ws.Range("E43") "  +0.08%  "
Set-TextValue $ws.Range("D44") "14.30"
Set-TextValue $ws.Range("E44") "  +1.19%  "
Set-TextValue $ws.Range("D45") "0.6420"
Set-TextValue $ws.Range("E45") "  -0.29%  "
Set-TextValue $ws.Range("D46") "2.293"
Set-TextValue $ws.Range("E46") "  -0.90%  "
Set-TextValue $ws.Range("B47") "PancakeSwap"
Set-TextValue $ws.Range("C47") "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws.Range("D47") "3.680"
Set-TextValue $ws.Range("E47") "  +0.07%  "
Set-TextValue $ws.Range("B48") "BabyDogeCoin"
Set-TextValue $ws.Range("C48") "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue $ws.Range("D48") "0.00000000361"
Set-TextValue $ws.Range("E48") "  +4.06%  "
Set-TextValue $ws.Range("B49") "Aave"
Set-TextValue $ws.Range("C49") "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Range("D49") "83.59"
Set-TextValue $ws.Range("E49") "  +0.65%  "
Set-TextValue $ws.Range("B50") "EOS"
Set-TextValue $ws.Range("C50") "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
Set-TextValue $ws.Range("D50") "1.225"
Set-TextValue $ws.Range("E50") "  -2.09%  "
Set-TextValue $ws.Range("D51") "0.07223"
Set-TextValue $ws.Range("E51") "  +0.57%  "
